# Verify Voltage Drop Calculation on adding devices in Multiple loops.xlsx
# - Refresh the recorded "Voltage Drop (V)" / "Voltage Drop (Worst Case)" values
#   for the first device row in Loop A and Loop B (re-run of the voltage-drop
#   calc after the "loading details" lookup/xpath rework).
# - Loop A's row 7 ("801 CH") drop value was previously typed as text ("0.70");
#   it is now a proper numeric value (0.3) like every other reading.
# - Active sheet/selection moved from Loop A to Loop B while reviewing results.

$wb = $excel.ActiveWorkbook

$wsA = $wb.Worksheets.Item("Add Devices Loop A")
$wsB = $wb.Worksheets.Item("Add Devices Loop B")

# --- Loop A ("Add Devices Loop A") ------------------------------------------
# Row 6 - LPBS 3000
$wsA.Range("F6").Value = 0.26
$wsA.Range("G6").Value = 0.47

# Row 7 - 801 CH (F7 was a text "0.70" value; now a real number)
$wsA.Range("F7").Value = 0.3
$wsA.Range("G7").Value = 0.55
# G7 picks up the "no bottom border" row style (matches the rest of the block)
$wsA.Range("G7").Borders.Item(9).LineStyle = -4142

# --- Loop B ("Add Devices Loop B") ------------------------------------------
# Row 6 - LPBS 3000
$wsB.Range("F6").Value = 0.5
$wsB.Range("G6").Value = 0.94

# --- Selection / active tab --------------------------------------------------
# Focus moves from Loop A to Loop B.
$wsA.Activate()
$wsA.Range("G5").Select() | Out-Null

$wsB.Activate()
$wsB.Range("F5").Select() | Out-Null
